$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values for rows 2-5 and 7 (row 6 unchanged at 0)
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 1
$ws.Range("G7").Value = 1
